{"js": "// Insert a new bulleted list item right after the paragraph that reads\n// \"(Familia, hermano, padre, tio);\" containing the new Empleo example,\n// matching the formatting (numbering / indentation) of the surrounding\n// list items.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is the \"(Familia, hermano, padre, tio);\" line.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"(Familia, hermano, padre, tio);\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"(Familia, hermano, padre, tio);\"');\n}\n\n// Insert the new list item immediately after it. insertParagraph(\"After\")\n// duplicates the paragraph-level formatting (numbering / indent / borders)\n// of the anchor paragraph, matching the rest of the bulleted list.\ntarget.insertParagraph(\n  \"(Empleo, (Pedro, trabajaPara, unEmpleador), (unEmpleador, brindaServiciosPara, unCliente), (Pedro, brindaServiciosPara, unCliente));\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item right after the paragraph that reads\n# \"(Familia, hermano, padre, tio);\" containing the new Empleo example,\n# matching the formatting (numbering / indentation) of the surrounding\n# list items.\n\n$d = $word.ActiveDocument\n\n$searchText = \"(Familia, hermano, padre, tio);\"\n$newText = \"(Empleo, (Pedro, trabajaPara, unEmpleador), (unEmpleador, brindaServiciosPara, unCliente), (Pedro, brindaServiciosPara, unCliente));\"\n\n# Locate the anchor paragraph using Find against the document's content range.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($searchText)\n\nif ($found) {\n    # $rng now covers the matched text; Paragraphs(1) is the paragraph that\n    # contains it.\n    $target = $rng.Paragraphs(1)\n\n    # Collapse to the end of the paragraph (includes its paragraph mark) and\n    # insert a brand-new paragraph after it. InsertParagraphAfter() on a\n    # Range duplicates the paragraph formatting (numbering, indents,\n    # borders, shading) of the paragraph it follows, same as pressing Enter\n    # at the end of the line in Word.\n    $insertionPoint = $target.Range\n    $insertionPoint.Collapse(0)\n    $insertionPoint.InsertParagraphAfter()\n\n    $newPara = $target.Next()\n    $newPara.Range.Text = $newText\n} else {\n    Write-Output \"WARNING: anchor paragraph not found\"\n}\n"}
